$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A61").Value = "GRT-USD"
